$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(0.1871, 0.1871, 0.1544000059366226, 0.3486000001430511, 0.0835999995470047)
    3 = @(0.2145, 0.2145, 0.2126, 0.3691000044345856, 0.03449999913573265)
    4 = @(0.2232, 0.2232, 0.2255, 0.3556999862194061, 0.09650000184774399)
    5 = @(0.2214, 0.2214, 0.2158, 0.4212000072002411, 0.1211000010371208)
    6 = @(0.2705, 0.2705, 0.2738, 0.4047999978065491, 0.18019999563694)
    7 = @(0.2618, 0.2618, 0.2655999958515167, 0.3228999972343445, 0.1555999964475632)
    8 = @(0.1858, 0.1858, 0.1848, 0.2567999958992004, 0.1035000011324883)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($i = 0; $i -lt $values.Length; $i++) {
        $col = 2 + $i
        $ws.Cells.Item($row, $col).Value = $values[$i]
    }
}
